$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move column D values (rows 2-4) down into column C as new rows 5-7.
$ws.Range("C5").Value2 = $ws.Range("D2").Value2
$ws.Range("C6").Value2 = $ws.Range("D3").Value2
$ws.Range("C7").Value2 = $ws.Range("D4").Value2

# Remove the now-redundant column D entirely (header + data).
$ws.Columns.Item(4).Delete()
